$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Command-line related words that have now been tested (column B holds the
# "Status" of each vocabulary word: "Coded" -> "Tested"). Two rows (the
# command-line related words "a>b b>a" / "ab>r", rows 31-32) are still only
# "Coded" - work on testing those from the command line has only just
# started, so they are deliberately left alone.
$skipRows = 31, 32

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    if ($skipRows -contains $r) {
        continue
    }

    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Text -eq "Coded") {
        $cell.Value = "Tested"
        $cell.WrapText = $false
    }
}

$ws.Range("B1").Select()
